# Automatic tracker update
# 1) Fill in "resultado" (G) / "profit" (H) for matches that were pending
#    and have now been settled.
# 2) Append the newly scraped fixtures as new rows at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Settle pending results --------------------------------------------
# Each entry: row, resultado, profit
$settled = @(
    @(417, "Fallo",   -1),
    @(419, "Acierto",  5),
    @(421, "Acierto",  1.5),
    @(425, "Fallo",   -1),
    @(426, "Fallo",   -1),
    @(428, "Acierto",  0.91),
    @(429, "Fallo",   -1),
    @(430, "Fallo",   -1),
    @(431, "Fallo",   -1),
    @(432, "Fallo",   -1),
    @(433, "Acierto",  0.91),
    @(434, "Acierto",  8),
    @(435, "Fallo",   -1),
    @(436, "Fallo",   -1),
    @(437, "Fallo",   -1),
    @(443, "Fallo",   -1)
)

foreach ($item in $settled) {
    $r = $item[0]
    $ws.Cells.Item($r, 7).Value = $item[1]
    $ws.Cells.Item($r, 8).Value = $item[2]
}

# --- 2) Append newly scraped fixtures --------------------------------------
# Each entry: row, event_id, fecha, jugador_A, jugador_B, pronostico, cuota
$newRows = @(
    @(445, 14540604, "2025-08-26", "Stefano Travaglia", "Gabriele Pennaforti", "Gana Gabriele Pennaforti", 3.4),
    @(446, 14528383, "2025-08-26", "Mili Poljičak", "Guy Den Ouden", "Gana Mili Poljičak", 2.63),
    @(447, 14528385, "2025-08-26", "Dmitry Popko", "Dimitar Kuzmanov", "Gana Dmitry Popko", 2.75),
    @(448, 14506232, "2025-08-26", "Marton Fucsovics", "Denis Shapovalov", "Gana Denis Shapovalov", 1.8),
    @(449, 14506225, "2025-08-26", "Quentin Halys", "David Goffin", "Gana David Goffin", 2.63),
    @(450, 14520005, "2025-08-26", "Leandro Riedi", "Pedro Martinez", "Gana Pedro Martinez", 2.2),
    @(451, 14506230, "2025-08-26", "Lorenzo Sonego", "Tristan Schoolkate", "Gana Tristan Schoolkate", 3.2),
    @(452, 14520011, "2025-08-26", "Billy Harris", "Felix Auger-Aliassime", "Gana Billy Harris", 8),
    @(453, 14520004, "2025-08-26", "Hugo Gaston", "Shintaro Mochizuki", "Gana Hugo Gaston", 2.5),
    @(454, 14506220, "2025-08-26", "Roberto Bautista Agut", "Jacob Fearnley", "Gana Jacob Fearnley", 2.1),
    @(455, 14506248, "2025-08-26", "Beatriz Haddad Maia", "Sonay Kartal", "Gana Beatriz Haddad Maia", 2.1),
    @(456, 14506242, "2025-08-26", "Caroline Dolehide", "Xinyu Wang", "Gana Xinyu Wang", 1.67),
    @(457, 14506246, "2025-08-26", "Lois Boisson", "Viktorija Golubic", "Gana Lois Boisson", 2.2),
    @(458, 14506261, "2025-08-26", "Sorana Cirstea", "Solana Sierra", "Gana Solana Sierra", 3.2),
    @(459, 14506253, "2025-08-26", "Donna Vekić", "Jessica Bouzas Maneiro", "Gana Donna Vekić", 2.2),
    @(460, 14506255, "2025-08-26", "Greet Minnen", "Naomi Osaka", "Gana Greet Minnen", 6.5),
    @(461, 14540602, "2025-08-26", "Federico Bondioli", "Giovanni Fonio", "Gana Federico Bondioli", 2.5),
    @(462, 14534338, "2025-08-26", "Borna Gojo", "Saba Purtseladze", "Gana Saba Purtseladze", 3),
    @(463, 14487518, "2025-08-27", "Alex Barrena", "Nicolas Alvarez Varona", "Gana Nicolas Alvarez Varona", 2.25),
    @(464, 14487522, "2025-08-27", "Gilles Arnaud Bailly", "Pedro Araujo", "Gana Pedro Araujo", 3.5),
    @(465, 14540525, "2025-08-26", "Adria Soriano Barrera", "Stuart Parker", "Gana Stuart Parker", 2.1),
    @(466, 14540526, "2025-08-26", "Jakub Paul", "Christian Langmo", "Gana Christian Langmo", 3)
)

foreach ($item in $newRows) {
    $r = $item[0]

    $ws.Cells.Item($r, 1).Value = $item[1]

    # The "fecha" column holds plain text dates (e.g. "2025-08-26"), not real
    # Excel dates. Force text storage so it isn't auto-converted to a date
    # serial number, then drop the number-format override so the cell keeps
    # the workbook's default (unstyled) look, matching the rest of the sheet.
    $dateCell = $ws.Cells.Item($r, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $item[2]
    $dateCell.ClearFormats()

    $ws.Cells.Item($r, 3).Value = $item[3]
    $ws.Cells.Item($r, 4).Value = $item[4]
    $ws.Cells.Item($r, 5).Value = $item[5]
    $ws.Cells.Item($r, 6).Value = $item[6]
}
